$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D keeps its original text representation (prices stored as
# plain text, not numbers) by forcing a text number-format before writing
# any numeric-looking strings, then clearing the format override again so
# the cells end up with no explicit style, matching the source workbook.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "39.948.53"
$ws.Range("E2").Value = "  +1.06%  "
$ws.Range("D3").Value = "2.225.56"
$ws.Range("E3").Value = "  +0.38%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "292.03"
$ws.Range("E5").Value = "  -2.02%  "
$ws.Range("D6").Value = "86.61"
$ws.Range("E6").Value = "  +4.24%  "
$ws.Range("E7").Value = "  +0.62%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "0.473"
$ws.Range("E9").Value = "  +0.76%  "
$ws.Range("D10").Value = "30.92"
$ws.Range("E10").Value = "  +5.50%  "
$ws.Range("D11").Value = "0.0786"
$ws.Range("E11").Value = "  +1.08%  "
$ws.Range("D12").Value = "47.16"
$ws.Range("E12").Value = "  -1.01%  "
$ws.Range("D13").Value = "0.108"
$ws.Range("E13").Value = "  +1.46%  "
$ws.Range("D14").Value = "6.35"
$ws.Range("E14").Value = "  +0.56%  "
$ws.Range("D15").Value = "2.574.02"
$ws.Range("E15").Value = "  +0.70%  "
$ws.Range("D16").Value = "14.13"
$ws.Range("E16").Value = "  -0.11%  "
$ws.Range("D17").Value = "2.229.48"
$ws.Range("E17").Value = "  +1.12%  "
$ws.Range("D18").Value = "0.729"
$ws.Range("E18").Value = "  +1.79%  "
$ws.Range("D19").Value = "39.923.22"
$ws.Range("E19").Value = "  +1.26%  "
$ws.Range("D20").Value = "0.0₃0885"
$ws.Range("E20").Value = "  +0.92%  "
$ws.Range("D21").Value = "11.15"
$ws.Range("E21").Value = "  +6.91%  "
$ws.Range("D22").Value = "5.82"
$ws.Range("E22").Value = "  +1.37%  "
$ws.Range("D23").Value = "65.75"
$ws.Range("E23").Value = "  +0.74%  "
$ws.Range("D24").Value = "236.52"
$ws.Range("E24").Value = "  +3.97%  "
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  -0.13%  "
$ws.Range("D26").Value = "2.47"
$ws.Range("E26").Value = "  +2.30%  "
$ws.Range("D27").Value = "1.84"
$ws.Range("E27").Value = "  +2.14%  "
$ws.Range("D28").Value = "22.82"
$ws.Range("E28").Value = "  +0.48%  "
$ws.Range("D29").Value = "2.21"
$ws.Range("E29").Value = "  +2.01%  "
$ws.Range("D30").Value = "9.27"
$ws.Range("E30").Value = "  +1.55%  "
$ws.Range("D31").Value = "33.07"
$ws.Range("E31").Value = "  +3.27%  "
$ws.Range("D32").Value = "151.46"
$ws.Range("E32").Value = "  +1.37%  "
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("D34").Value = "4.95"
$ws.Range("E34").Value = "  +1.79%  "
$ws.Range("D35").Value = "0.0719"
$ws.Range("E35").Value = "  +3.07%  "
$ws.Range("D36").Value = "2.37"
$ws.Range("E36").Value = "  +1.76%  "
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").Value = "2.82"
$ws.Range("E37").Value = "  +6.94%  "
$ws.Range("B38").Value = "Celestia"
$ws.Range("C38").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D38").Value = "16.06"
$ws.Range("E38").Value = "  +3.53%  "
$ws.Range("D39").Value = "0.111"
$ws.Range("E39").Value = "  +1.37%  "
$ws.Range("D40").Value = "0.0998"
$ws.Range("E40").Value = "  +2.84%  "
$ws.Range("D41").Value = "1.70"
$ws.Range("E41").Value = "  +3.00%  "
$ws.Range("D42").Value = "3.80"
$ws.Range("E42").Value = "  +3.94%  "
$ws.Range("D43").Value = "2.063.87"
$ws.Range("E43").Value = "  +8.53%  "
$ws.Range("D44").Value = "18.57"
$ws.Range("E44").Value = "  +14.08%  "
$ws.Range("E45").Value = "  +3.18%  "
$ws.Range("D46").Value = "0.0268"
$ws.Range("E46").Value = "  +3.24%  "
$ws.Range("D47").Value = "9.93"
$ws.Range("E47").Value = "  +9.99%  "
$ws.Range("D48").Value = "2.60"
$ws.Range("E48").Value = "  -1.51%  "
$ws.Range("D49").Value = "2.451.61"
$ws.Range("E49").Value = "  +1.26%  "
$ws.Range("D50").Value = "72.31"
$ws.Range("E50").Value = "  +2.01%  "
$ws.Range("D51").Value = "89.25"
$ws.Range("E51").Value = "  +2.04%  "

$priceRange.Style = "Normal"

